# Updates cryptos list values (price/volume columns, plus the Monero/Filecoin row swap)
# to match the GitHub-Actions-refreshed data. Every target cell holds a plain text value
# (col D prices, col E volumes use a "  +x.xx%  " padded format) - a leading single-quote
# forces Excel to keep the assignment as text instead of auto-parsing it as a number,
# then resetting the Style back to "Normal" drops the quote-prefix formatting again so
# no stray cell style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "42.880.93"
Set-TextValue "E2" "  +2.78%  "
Set-TextValue "D3" "2.304.03"
Set-TextValue "E3" "  +0.89%  "
Set-TextValue "E4" "  +0.13%  "
Set-TextValue "D5" "319.12"
Set-TextValue "E5" "  +1.47%  "
Set-TextValue "D6" "104.68"
Set-TextValue "E6" "  +2.20%  "
Set-TextValue "E7" "  +0.79%  "
Set-TextValue "E8" "  +0.11%  "
Set-TextValue "E9" "  +1.62%  "
Set-TextValue "D10" "40.19"
Set-TextValue "E10" "  +3.69%  "
Set-TextValue "D11" "0.0910"
Set-TextValue "E11" "  +0.73%  "
Set-TextValue "D12" "8.64"
Set-TextValue "E12" "  +5.05%  "
Set-TextValue "E13" "  +1.30%  "
Set-TextValue "D14" "0.977"
Set-TextValue "E14" "  +1.98%  "
Set-TextValue "D15" "15.39"
Set-TextValue "E15" "  +1.07%  "
Set-TextValue "D16" "2.652.83"
Set-TextValue "E16" "  +0.96%  "
Set-TextValue "D17" "2.301.91"
Set-TextValue "E17" "  +0.85%  "
Set-TextValue "D18" "42.799.21"
Set-TextValue "E18" "  +2.65%  "
Set-TextValue "D19" "7.54"
Set-TextValue "E19" "  +2.01%  "
Set-TextValue "D20" "0.0000106"
Set-TextValue "E20" "  +1.30%  "
Set-TextValue "D21" "13.55"
Set-TextValue "E21" "  +33.97%  "
Set-TextValue "E22" "  +1.14%  "
Set-TextValue "D23" "3.56"
Set-TextValue "E23" "  -1.57%  "
Set-TextValue "D24" "271.14"
Set-TextValue "E24" "  -2.92%  "
Set-TextValue "E25" "  +0.83%  "
Set-TextValue "E26" "  -0.53%  "
Set-TextValue "E27" "  +2.36%  "
Set-TextValue "E28" "  -2.32%  "
Set-TextValue "D29" "22.71"
Set-TextValue "E29" "  -0.68%  "
Set-TextValue "D30" "38.07"
Set-TextValue "E30" "  +9.18%  "
Set-TextValue "B31" "Filecoin"
Set-TextValue "C31" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "6.26"
Set-TextValue "E31" "  +8.10%  "
Set-TextValue "B32" "Monero"
Set-TextValue "C32" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D32" "165.78"
Set-TextValue "E32" "  +1.92%  "
Set-TextValue "D33" "0.0893"
Set-TextValue "E33" "  +3.01%  "
Set-TextValue "D34" "0.133"
Set-TextValue "E34" "  +1.27%  "
Set-TextValue "E35" "  +0.43%  "
Set-TextValue "D36" "2.55"
Set-TextValue "E36" "  -11.29%  "
Set-TextValue "D37" "4.63"
Set-TextValue "E37" "  +2.42%  "
Set-TextValue "E38" "  +2.97%  "
Set-TextValue "E39" "  +2.66%  "
Set-TextValue "E40" "  -4.73%  "
Set-TextValue "D41" "1.57"
Set-TextValue "E41" "  +8.14%  "
Set-TextValue "D42" "100.45"
Set-TextValue "E42" "  +0.60%  "
Set-TextValue "D43" "70.71"
Set-TextValue "E43" "  +2.00%  "
Set-TextValue "D44" "0.227"
Set-TextValue "E44" "  +1.73%  "
Set-TextValue "E45" "  +0.19%  "
Set-TextValue "D46" "12.36"
Set-TextValue "E46" "  +4.43%  "
Set-TextValue "D47" "82.93"
Set-TextValue "E47" "  +9.76%  "
Set-TextValue "E48" "  -0.77%  "
Set-TextValue "D49" "5.32"
Set-TextValue "E49" "  +1.76%  "
Set-TextValue "E50" "  -0.22%  "
Set-TextValue "D51" "1.596.36"
Set-TextValue "E51" "  +4.81%  "
